$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'1582"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "'3870086.96"
$ws.Range("D2").Style = "Normal"

$ws.Range("C4").Value = "'1160"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'4867950.40"
$ws.Range("D4").Style = "Normal"

$ws.Range("C6").Value = "'795"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'3096362.16"
$ws.Range("D6").Style = "Normal"

$ws.Range("C8").Value = "'42"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'192144.45"
$ws.Range("D8").Style = "Normal"

$ws.Range("C9").Value = "'243"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'780843.05"
$ws.Range("D9").Style = "Normal"

$ws.Range("C10").Value = "'4"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'17000.00"
$ws.Range("D10").Style = "Normal"

$ws.Range("C11").Value = "'429"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'2246809.45"
$ws.Range("D11").Style = "Normal"

$ws.Range("C12").Value = "'203"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'1177743.34"
$ws.Range("D12").Style = "Normal"

$ws.Range("C15").Value = "'264"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'814024.01"
$ws.Range("D15").Style = "Normal"

$ws.Range("C17").Value = "'580"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3446481.24"
$ws.Range("D17").Style = "Normal"

$ws.Range("C18").Value = "'166"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'767616.42"
$ws.Range("D18").Style = "Normal"

$ws.Range("C19").Value = "'17"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'42471.00"
$ws.Range("D19").Style = "Normal"

$ws.Range("C20").Value = "'10"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'40418.77"
$ws.Range("D20").Style = "Normal"

$ws.Range("C21").Value = "'220"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'616245.00"
$ws.Range("D21").Style = "Normal"

$ws.Range("C23").Value = "'411"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'2265742.61"
$ws.Range("D23").Style = "Normal"

$ws.Range("C24").Value = "'187"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'865657.27"
$ws.Range("D24").Style = "Normal"

$ws.Range("C31").Value = "'392"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'1106179.09"
$ws.Range("D31").Style = "Normal"

$ws.Range("C32").Value = "'9"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'40000.00"
$ws.Range("D32").Style = "Normal"

$ws.Range("C33").Value = "'742"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'4609230.94"
$ws.Range("D33").Style = "Normal"

$ws.Range("C35").Value = "'493"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'2539661.28"
$ws.Range("D35").Style = "Normal"

$ws.Range("C37").Value = "'20"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'65432.00"
$ws.Range("D37").Style = "Normal"

$ws.Range("C38").Value = "'536"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'1395802.65"
$ws.Range("D38").Style = "Normal"

$ws.Range("C39").Value = "'252"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'861292.04"
$ws.Range("D39").Style = "Normal"

$ws.Range("C40").Value = "'253"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'839460.95"
$ws.Range("D40").Style = "Normal"

$ws.Range("C43").Value = "'339"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'1110647.40"
$ws.Range("D43").Style = "Normal"

$ws.Range("C44").Value = "'149"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'842560.67"
$ws.Range("D44").Style = "Normal"

$ws.Range("C45").Value = "'222"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'1015692.25"
$ws.Range("D45").Style = "Normal"

$ws.Range("C48").Value = "'641"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'1844898.85"
$ws.Range("D48").Style = "Normal"

$ws.Range("C50").Value = "'913"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'5496216.43"
$ws.Range("D50").Style = "Normal"

$ws.Range("C51").Value = "'645"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'3391902.84"
$ws.Range("D51").Style = "Normal"

$ws.Range("C53").Value = "'32"
$ws.Range("C53").Style = "Normal"
$ws.Range("D53").Value = "'177011.07"
$ws.Range("D53").Style = "Normal"

$ws.Range("C54").Value = "'8342"
$ws.Range("C54").Style = "Normal"
$ws.Range("D54").Value = "'21422749.15"
$ws.Range("D54").Style = "Normal"

$ws.Range("C57").Value = "'39"
$ws.Range("C57").Style = "Normal"
$ws.Range("D57").Value = "'277600.00"
$ws.Range("D57").Style = "Normal"

$ws.Range("C58").Value = "'5703"
$ws.Range("C58").Style = "Normal"
$ws.Range("D58").Value = "'25226744.02"
$ws.Range("D58").Style = "Normal"

$ws.Range("C59").Value = "'15"
$ws.Range("C59").Style = "Normal"
$ws.Range("D59").Value = "'121500.00"
$ws.Range("D59").Style = "Normal"

$ws.Range("C60").Value = "'5669"
$ws.Range("C60").Style = "Normal"
$ws.Range("D60").Value = "'21450984.74"
$ws.Range("D60").Style = "Normal"

$ws.Range("C61").Value = "'61"
$ws.Range("C61").Style = "Normal"
$ws.Range("D61").Value = "'166670.00"
$ws.Range("D61").Style = "Normal"

$ws.Range("C62").Value = "'117"
$ws.Range("C62").Style = "Normal"
$ws.Range("D62").Value = "'489315.40"
$ws.Range("D62").Style = "Normal"

$ws.Range("C74").Value = "'268"
$ws.Range("C74").Style = "Normal"
$ws.Range("D74").Value = "'812360.35"
$ws.Range("D74").Style = "Normal"

$ws.Range("C75").Value = "'461"
$ws.Range("C75").Style = "Normal"
$ws.Range("D75").Value = "'2437807.65"
$ws.Range("D75").Style = "Normal"

$ws.Range("C76").Value = "'272"
$ws.Range("C76").Style = "Normal"
$ws.Range("D76").Value = "'1612534.51"
$ws.Range("D76").Style = "Normal"

$ws.Range("C79").Value = "'434"
$ws.Range("C79").Style = "Normal"
$ws.Range("D79").Value = "'1255768.80"
$ws.Range("D79").Style = "Normal"

$ws.Range("C81").Value = "'1139"
$ws.Range("C81").Style = "Normal"
$ws.Range("D81").Value = "'6972462.63"
$ws.Range("D81").Style = "Normal"

$ws.Range("C82").Value = "'607"
$ws.Range("C82").Style = "Normal"
$ws.Range("D82").Value = "'3186176.81"
$ws.Range("D82").Style = "Normal"

$ws.Range("C83").Value = "'46"
$ws.Range("C83").Style = "Normal"
$ws.Range("D83").Value = "'172078.00"
$ws.Range("D83").Style = "Normal"

$ws.Range("C84").Value = "'41"
$ws.Range("C84").Style = "Normal"
$ws.Range("D84").Value = "'182727.77"
$ws.Range("D84").Style = "Normal"

$ws.Range("C85").Value = "'742"
$ws.Range("C85").Style = "Normal"
$ws.Range("D85").Value = "'1862718.91"
$ws.Range("D85").Style = "Normal"

$ws.Range("C88").Value = "'1067"
$ws.Range("C88").Style = "Normal"
$ws.Range("D88").Value = "'4538945.84"
$ws.Range("D88").Style = "Normal"

$ws.Range("C89").Value = "'762"
$ws.Range("C89").Style = "Normal"
$ws.Range("D89").Value = "'2926330.18"
$ws.Range("D89").Style = "Normal"

$ws.Range("C91").Value = "'35"
$ws.Range("C91").Style = "Normal"
$ws.Range("D91").Value = "'123571.23"
$ws.Range("D91").Style = "Normal"

$ws.Range("C92").Value = "'247"
$ws.Range("C92").Style = "Normal"
$ws.Range("D92").Value = "'583350.00"
$ws.Range("D92").Style = "Normal"

$ws.Range("C94").Value = "'580"
$ws.Range("C94").Style = "Normal"
$ws.Range("D94").Value = "'2774601.99"
$ws.Range("D94").Style = "Normal"

$ws.Range("C95").Value = "'211"
$ws.Range("C95").Style = "Normal"
$ws.Range("D95").Value = "'796666.11"
$ws.Range("D95").Style = "Normal"

$ws.Range("C96").Value = "'18"
$ws.Range("C96").Style = "Normal"
$ws.Range("D96").Value = "'60500.00"
$ws.Range("D96").Style = "Normal"

$ws.Range("C97").Value = "'9"
$ws.Range("C97").Style = "Normal"
$ws.Range("D97").Value = "'34670.00"
$ws.Range("D97").Style = "Normal"

$ws.Range("C98").Value = "'1238"
$ws.Range("C98").Style = "Normal"
$ws.Range("D98").Value = "'3121613.28"
$ws.Range("D98").Style = "Normal"

$ws.Range("C100").Value = "'11"
$ws.Range("C100").Style = "Normal"
$ws.Range("D100").Value = "'37560.00"
$ws.Range("D100").Style = "Normal"

$ws.Range("C101").Value = "'1422"
$ws.Range("C101").Style = "Normal"
$ws.Range("D101").Value = "'6454511.94"
$ws.Range("D101").Style = "Normal"

$ws.Range("C103").Value = "'1384"
$ws.Range("C103").Style = "Normal"
$ws.Range("D103").Value = "'5715262.90"
$ws.Range("D103").Style = "Normal"

$ws.Range("C105").Value = "'76"
$ws.Range("C105").Style = "Normal"
$ws.Range("D105").Value = "'297729.61"
$ws.Range("D105").Style = "Normal"
